# Generate Report for Handoff
#
# Updates the "latest handoff" timestamps for the
# 66f8ed3e-b652-4472-84d9-8a2a16539b78 file across the Overview sheet and
# each per-language handoff-status sheet, as produced by a fresh report
# generation run.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 is the 66f8ed3e-... file; column D is
# "Latest Handoff Date".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-27-11 18:27:06"

# zh-cn sheet: row 7 is the 66f8ed3e-... file; column E is
# "Latest Handoff Datetime".
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-11 18:27:03"

# de-de sheet: row 7 is the 66f8ed3e-... file; column E is
# "Latest Handoff Datetime".
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-11 18:27:06"
